# This script reproduces the commit 'culture collection removed from MIxS, per
# INSDC2017 review': the 'culture_collection' header/column (AP) is deleted,
# shifting every later column (and its shared string) one slot to the left.
# This runtime does not carry per-cell header comments along with a column
# delete, so the comments are re-applied explicitly afterwards to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of (cell, text) pairs: after column AP is removed, the comment
# that used to sit one column to the right now belongs at this cell.
$shiftedComments = @(
    @('AP15', 'density of sample'),
    @('AQ15', 'concentration of diether lipids; can include multiple types of diether lipids'),
    @('AR15', 'concentration of dissolved carbon dioxide'),
    @('AS15', 'concentration of dissolved hydrogen'),
    @('AT15', 'dissolved inorganic carbon concentration'),
    @('AU15', 'concentration of dissolved inorganic nitrogen'),
    @('AV15', 'concentration of dissolved inorganic phosphorus'),
    @('AW15', 'concentration of dissolved organic carbon'),
    @('AX15', 'dissolved organic nitrogen concentration measured as; total dissolved nitrogen - NH4 - NO3 - NO2'),
    @('AY15', 'concentration of dissolved oxygen'),
    @('AZ15', 'visible waveband radiance and irradiance measurements in the water column'),
    @('BA15', 'The elevation of the sampling site as measured by the vertical distance from mean sea level.'),
    @('BB15', 'Plasmids that have significance phenotypic consequence'),
    @('BC15', 'raw or converted fluorescence of water'),
    @('BD15', 'measurement of glucosidase activity'),
    @('BE15', 'Health or disease status of sample at time of collection'),
    @('BF15', 'The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, "Homo sapiens".'),
    @('BG15', 'NCBI taxonomy ID of the host, e.g. 9606'),
    @('BH15', 'Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.'),
    @('BI15', 'measurement of light intensity'),
    @('BJ15', 'A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html'),
    @('BK15', 'concentration of magnesium'),
    @('BL15', 'measurement of mean friction velocity'),
    @('BM15', 'measurement of mean peak friction velocity'),
    @('BN15', 'any other measurement performed or parameter collected, that is not listed here'),
    @('BO15', 'concentration of n-alkanes; can include multiple n-alkanes'),
    @('BP15', 'concentration of nitrate'),
    @('BQ15', 'concentration of nitrite'),
    @('BR15', 'concentration of nitrogen (total)'),
    @('BS15', 'concentration of organic carbon'),
    @('BT15', 'concentration of organic matter'),
    @('BU15', 'concentration of organic nitrogen'),
    @('BV15', 'total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts'),
    @('BW15', 'oxygenation status of sample'),
    @('BX15', 'concentration of particulate organic carbon'),
    @('BY15', 'concentration of particulate organic nitrogen'),
    @('BZ15', 'To what is the entity pathogenic'),
    @('CA15', 'type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types'),
    @('CB15', 'concentration of petroleum hydrocarbon'),
    @('CC15', 'pH measurement'),
    @('CD15', 'concentration of phaeopigments; can include multiple phaeopigments'),
    @('CE15', 'concentration of phosphate'),
    @('CF15', 'concentration of phospholipid fatty acids; can include multiple values'),
    @('CG15', 'measurement of photon flux'),
    @('CH15', 'concentration of potassium'),
    @('CI15', 'pressure to which the sample is subject, in atmospheres'),
    @('CJ15', 'measurement of primary production'),
    @('CK15', 'redox potential, measured relative to a hydrogen cell, indicating oxidation or reduction potential'),
    @('CL15', 'Primary publication or genome report in the form of pubmed ID, DOI or URL'),
    @('CM15', 'salinity measurement'),
    @('CN15', 'Method or device employed for collecting sample'),
    @('CO15', 'Processing applied to the sample during or after isolation'),
    @('CP15', 'Amount or size of sample (volume, mass or area) that was collected'),
    @('CQ15', 'duration for which sample was stored'),
    @('CR15', 'location at which sample was stored, usually name of a specific freezer/room'),
    @('CS15', 'temperature at which sample was stored, e.g. -80'),
    @('CT15', 'volume (mL) or weight (g) of sample processed for DNA extraction'),
    @('CU15', 'concentration of silicate'),
    @('CV15', 'sodium concentration'),
    @('CW15', 'concentration of soluble reactive phosphorus'),
    @('CX15', 'unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.'),
    @('CY15', 'Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier'),
    @('CZ15', 'Information about the genetic distinctness of the lineage (eg., biovar, serovar)'),
    @('DA15', 'concentration of sulfate'),
    @('DB15', 'concentration of sulfide'),
    @('DC15', 'concentration of suspended particulate matter'),
    @('DD15', 'temperature of the sample at time of sampling'),
    @('DE15', 'stage of tide'),
    @('DF15', 'measurement of total depth of water column'),
    @('DG15', 'total dissolved nitrogen concentration, reported as nitrogen, measured by: total dissolved nitrogen = NH4 + NO3NO2 + dissolved organic nitrogen'),
    @('DH15', 'total inorganic nitrogen content'),
    @('DI15', 'total nitrogen content of the sample'),
    @('DJ15', 'total particulate carbon content'),
    @('DK15', 'total phosphorus concentration, calculated by: total phosphorus = total dissolved phosphorus + particulate phosphorus. Can also be measured without filtering, reported as phosphorus'),
    @('DL15', 'Feeding position in food chain (eg., chemolithotroph)'),
    @('DM15', 'measurement of magnitude and direction of flow within a fluid')
)

# Step 1: delete the "culture_collection" column (AP), shifting later columns left
$ws.Columns("AP").Delete()

# Step 2: re-apply the shifted comment text to each affected header cell
foreach ($pair in $shiftedComments) {
    $cellRef = $pair[0]
    $commentText = $pair[1]
    $ws.Range($cellRef).Comment.Text($commentText) | Out-Null
}

# Step 3: drop the now-orphaned comment that used to belong to the last column
$ws.Range('DN15').Comment.Delete()

Write-Host "culture_collection column removed. UsedRange:" $ws.UsedRange.Address()
